$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 2666.3333
$ws.Range("I69").Value = 2000
$ws.Range("J69").Value = 3999
$ws.Range("K69").Value = 6000
$ws.Range("L69").Value = 11997
$ws.Range("M69").Value = -5126
$ws.Range("N69").Value = -13745
$ws.Range("H72").Value = 2666.3333
$ws.Range("I72").Value = 2000
$ws.Range("J72").Value = 3999
$ws.Range("K72").Value = 18000
$ws.Range("L72").Value = 35991
$ws.Range("M72").Value = -13632
$ws.Range("N72").Value = -44727
$ws.Range("H76").Value = 45458084
$ws.Range("I76").Value = 71431944
$ws.Range("J76").Value = 3825
$ws.Range("K76").Value = 71431944
$ws.Range("L76").Value = 3825
$ws.Range("M76").Value = -71431629
$ws.Range("N76").Value = -4455
$ws.Range("H79").Value = 45458084
$ws.Range("I79").Value = 71431944
$ws.Range("J79").Value = 3825
$ws.Range("K79").Value = 71431944
$ws.Range("L79").Value = 3825
$ws.Range("M79").Value = -71430852
$ws.Range("N79").Value = -6009
$ws.Range("H107").Value = 139.15384
$ws.Range("I107").Value = 152.6
$ws.Range("K107").Value = 152.6
$ws.Range("M107").Value = 1767.4
$ws.Range("H117").Value = 94876.78
$ws.Range("J117").Value = 94876.78
$ws.Range("L117").Value = 94876.78
$ws.Range("N117").Value = -104054.78
$ws.Range("H123").Value = 53749.125
$ws.Range("J123").Value = 53749.125
$ws.Range("L123").Value = 53749.125
$ws.Range("N123").Value = -63549.125
$ws.Range("H134").Value = 99995
$ws.Range("J134").Value = 99995
$ws.Range("L134").Value = 99995
$ws.Range("N134").Value = -110135
$ws.Range("H136").Value = 99991
$ws.Range("J136").Value = 99991
$ws.Range("L136").Value = 99991
$ws.Range("N136").Value = -110191
$ws.Range("H138").Value = 1384.7188
$ws.Range("J138").Value = 2461
$ws.Range("L138").Value = 7383
$ws.Range("N138").Value = -17663

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 679.9091
$ws.Range("I2").Value = 572.375
$ws.Range("K2").Value = 572.375
$ws.Range("M2").Value = -459.375
$ws.Range("H7").Value = 89497
$ws.Range("J7").Value = 89497
$ws.Range("L7").Value = 89497
$ws.Range("N7").Value = -89725
$ws.Range("H32").Value = 5036.9
$ws.Range("I32").Value = 2025.2241
$ws.Range("J32").Value = 19593.334
$ws.Range("K32").Value = 2025.2241
$ws.Range("L32").Value = 19593.334
$ws.Range("M32").Value = -1738.2241
$ws.Range("N32").Value = -20167.334
$ws.Range("H45").Value = 20836560
$ws.Range("I45").Value = 3940.25
$ws.Range("J45").Value = 62501800
$ws.Range("K45").Value = 3940.25
$ws.Range("L45").Value = 62501800
$ws.Range("M45").Value = -3563.25
$ws.Range("N45").Value = -62502554
$ws.Range("H97").Value = 1070.875
$ws.Range("I97").Value = 938.2857
$ws.Range("K97").Value = 938.2857
$ws.Range("M97").Value = -442.2857
$ws.Range("H116").Value = 679.9091
$ws.Range("I116").Value = 572.375
$ws.Range("K116").Value = 572.375
$ws.Range("M116").Value = 1721.625
$ws.Range("H117").Value = 42944
$ws.Range("J117").Value = 42944
$ws.Range("L117").Value = 42944
$ws.Range("N117").Value = -52122

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 679.9091
$ws.Range("I3").Value = 572.375
$ws.Range("K3").Value = 572.375
$ws.Range("M3").Value = -458.375
$ws.Range("H52").Value = 99988
$ws.Range("J52").Value = 99988
$ws.Range("L52").Value = 99988
$ws.Range("N52").Value = -100514
$ws.Range("H107").Value = 2024.5
$ws.Range("J107").Value = 4877.5
$ws.Range("L107").Value = 4877.5
$ws.Range("N107").Value = -8717.5
$ws.Range("H115").Value = 70854.42999999999
$ws.Range("J115").Value = 72330
$ws.Range("L115").Value = 72330
$ws.Range("N115").Value = -75464
$ws.Range("H121").Value = 99988
$ws.Range("J121").Value = 99988
$ws.Range("L121").Value = 99988
$ws.Range("N121").Value = -103482
$ws.Range("H127").Value = 63293.332
$ws.Range("J127").Value = 63293.332
$ws.Range("L127").Value = 63293.332
$ws.Range("N127").Value = -73213.33199999999
$ws.Range("H135").Value = 45107.617
$ws.Range("J135").Value = 45107.617
$ws.Range("L135").Value = 45107.617
$ws.Range("N135").Value = -55247.617
$ws.Range("H138").Value = 99995
$ws.Range("J138").Value = 99995
$ws.Range("L138").Value = 99995
$ws.Range("N138").Value = -110275

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 32298.4
$ws.Range("J9").Value = 32298.4
$ws.Range("L9").Value = 32298.4
$ws.Range("N9").Value = -32634.4
$ws.Range("H107").Value = 650.1539
$ws.Range("I107").Value = 608.0909
$ws.Range("J107").Value = 881.5
$ws.Range("K107").Value = 608.0909
$ws.Range("L107").Value = 881.5
$ws.Range("M107").Value = 1311.9091
$ws.Range("N107").Value = -4721.5
$ws.Range("H117").Value = 39581.832
$ws.Range("J117").Value = 39581.832
$ws.Range("L117").Value = 39581.832
$ws.Range("N117").Value = -48759.832
$ws.Range("H122").Value = 2193.9546
$ws.Range("I122").Value = 1875.4667
$ws.Range("K122").Value = 5626.4001
$ws.Range("M122").Value = -3176.4001
$ws.Range("H138").Value = 70615.336
$ws.Range("J138").Value = 70615.336
$ws.Range("L138").Value = 70615.336
$ws.Range("N138").Value = -80895.336

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 3343399.8
$ws.Range("I9").Value = 3343399.8
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 10030199.4
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -10029975.4
$ws.Range("N9").ClearContents()
$ws.Range("H18").Value = 10109.272
$ws.Range("I18").Value = 11911.333
$ws.Range("K18").Value = 35733.999
$ws.Range("M18").Value = -35564.999
$ws.Range("H93").Value = 3333.3333
$ws.Range("I93").Value = 1250
$ws.Range("K93").Value = 3750
$ws.Range("M93").Value = -1878
$ws.Range("H109").Value = 2364.0588
$ws.Range("I109").Value = 1218.9
$ws.Range("K109").Value = 3656.7
$ws.Range("M109").Value = -2616.7
$ws.Range("H115").Value = 1050
$ws.Range("I115").Value = 1050
$ws.Range("K115").Value = 3150
$ws.Range("M115").Value = -1975
$ws.Range("H117").Value = 245.4
$ws.Range("J117").Value = 144.83333
$ws.Range("L117").Value = 434.49999
$ws.Range("N117").Value = -7318.49999

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 71431290
$ws.Range("I80").Value = 166668670
$ws.Range("J80").Value = 3249.375
$ws.Range("K80").Value = 166668670
$ws.Range("L80").Value = 3249.375
$ws.Range("M80").Value = -166667672
$ws.Range("N80").Value = -5245.375
$ws.Range("H83").Value = 71431290
$ws.Range("I83").Value = 166668670
$ws.Range("J83").Value = 3249.375
$ws.Range("K83").Value = 833343350
$ws.Range("L83").Value = 16246.875
$ws.Range("M83").Value = -833338358
$ws.Range("N83").Value = -26230.875
$ws.Range("H107").Value = 693.76666
$ws.Range("I107").Value = 619.36365
$ws.Range("J107").Value = 736.8421
$ws.Range("K107").Value = 619.36365
$ws.Range("L107").Value = 736.8421
$ws.Range("M107").Value = 1300.63635
$ws.Range("N107").Value = -4576.8421
$ws.Range("H108").Value = 51801.168
$ws.Range("J108").Value = 51801.168
$ws.Range("L108").Value = 51801.168
$ws.Range("N108").Value = -59481.168
$ws.Range("H109").Value = 27688.727
$ws.Range("J109").Value = 27688.727
$ws.Range("L109").Value = 27688.727
$ws.Range("N109").Value = -29768.727
$ws.Range("H110").Value = 68329.875
$ws.Range("J110").Value = 68329.875
$ws.Range("L110").Value = 68329.875
$ws.Range("N110").Value = -76509.875
$ws.Range("H116").Value = 59997.332
$ws.Range("J116").Value = 59997.332
$ws.Range("L116").Value = 59997.332
$ws.Range("N116").Value = -69175.33199999999
$ws.Range("H119").Value = 52109.777
$ws.Range("J119").Value = 52109.777
$ws.Range("L119").Value = 52109.777
$ws.Range("N119").Value = -61785.777
$ws.Range("H140").Value = 95748
$ws.Range("J140").Value = 97497.336
$ws.Range("L140").Value = 97497.336
$ws.Range("N140").Value = -107857.336

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3215.4285
$ws.Range("I46").Value = 2242.4285
$ws.Range("J46").Value = 3539.762
$ws.Range("K46").Value = 2242.4285
$ws.Range("L46").Value = 3539.762
$ws.Range("M46").Value = -2054.4285
$ws.Range("N46").Value = -3915.762
$ws.Range("H117").Value = 74940
$ws.Range("J117").Value = 79096.664
$ws.Range("L117").Value = 79096.664
$ws.Range("N117").Value = -88274.664
$ws.Range("H123").Value = 75152.664
$ws.Range("J123").Value = 78248
$ws.Range("L123").Value = 78248
$ws.Range("N123").Value = -88048
$ws.Range("H136").Value = 3417.875
$ws.Range("I136").Value = 4630.4443
$ws.Range("K136").Value = 13891.3329
$ws.Range("M136").Value = -11341.3329

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3166.1667
$ws.Range("I81").Value = 2749.5
$ws.Range("K81").Value = 5499
$ws.Range("M81").Value = -4438
$ws.Range("H84").Value = 3166.1667
$ws.Range("I84").Value = 2749.5
$ws.Range("K84").Value = 27495
$ws.Range("M84").Value = -22191
$ws.Range("H121").Value = 37354.855
$ws.Range("J121").Value = 37354.855
$ws.Range("L121").Value = 37354.855
$ws.Range("N121").Value = -40848.855
